$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.794.81"

$ws.Range("E2").Value = "  +2.10%  "

$ws.Range("D3").Value = "3.040.39"

$ws.Range("E3").Value = "  +1.81%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "523.04"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = "  +5.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.91"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = "  +5.32%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.442"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = "  +3.54%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.61"
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = "  +4.63%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.110"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = "  +5.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.364"
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = "  +2.88%  "

$ws.Range("E12").Value = "  +2.29%  "

$ws.Range("D13").Value = "3.550.36"

$ws.Range("E13").Value = "  +1.35%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.23"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = "  +4.84%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000163"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = "  +11.08%  "

$ws.Range("D16").Value = "57.655.08"

$ws.Range("E16").Value = "  +2.05%  "

$ws.Range("D17").Value = "3.065.02"

$ws.Range("E17").Value = "  +2.62%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.12"
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = "  +3.92%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.83"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = "  +3.73%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.07"
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = "  +4.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "336.48"
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = "  +3.53%  "

$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.494"
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = "  +6.34%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.84"
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = "  +5.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.174"
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = "  +5.47%  "

$ws.Range("B26").Value = "PEPE"

$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"

$ws.Range("D26").Value = "0.0₃0941"

$ws.Range("E26").Value = "  +5.64%  "

$ws.Range("B27").Value = "Binance-PegBSC-USD"

$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.989"
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = "  -0.73%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.89"
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = "  +3.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.21"
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = "  +5.61%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.84"
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = "  +6.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.23"
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = "  +2.72%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.91"
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = "  +4.87%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "158.20"
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = "  +2.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.70"
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = "  +3.97%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.83"
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = "  +3.19%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.31"
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = "  +1.30%  "

$ws.Range("B37").Value = "EnergySwap"

$ws.Range("C37").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "24.83"
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = "  +6.24%  "

$ws.Range("B38").Value = "Hedera"

$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0688"
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = "  +2.52%  "

$ws.Range("D39").Value = "3.062.50"

$ws.Range("E39").Value = "  +1.38%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.61"
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = "  +0.86%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.997"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = "  -0.25%  "

$ws.Range("D42").Value = "2.332.48"

$ws.Range("E42").Value = "  +5.92%  "

$ws.Range("B43").Value = "Filecoin"

$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.79"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = "  +6.10%  "

$ws.Range("B44").Value = "Mantle"

$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.657"
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = "  +2.69%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.46"
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = "  +2.54%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.02"
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = "  +0.36%  "

$ws.Range("E47").Value = "  +4.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0245"
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = "  +3.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.96"
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = "  +5.77%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.81"
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = "  +2.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0891"
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = "  +4.58%  "
